# Update "想去人数" (want-to-go count) values in column F across all sheets,
# reflecting a refreshed data scrape (gh-pages output regenerated).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 16
$ws.Range("F9").Value = 726
$ws.Range("F10").Value = 2675
$ws.Range("F11").Value = 2675
$ws.Range("F13").Value = 1754
$ws.Range("F14").Value = 607
$ws.Range("F15").Value = 273
$ws.Range("F16").Value = 686
$ws.Range("F17").Value = 5030
$ws.Range("F18").Value = 208
$ws.Range("F21").Value = 3385
$ws.Range("F25").Value = 37
$ws.Range("F26").Value = 2417
$ws.Range("F28").Value = 366
$ws.Range("F31").Value = 480
$ws.Range("F32").Value = 1295
$ws.Range("F34").Value = 5
$ws.Range("F37").Value = 54
$ws.Range("F38").Value = 1430
$ws.Range("F39").Value = 11
$ws.Range("F40").Value = 1386
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 137
$ws.Range("F13").Value = 77
$ws.Range("F16").Value = 139
$ws.Range("F17").Value = 326
$ws.Range("F18").Value = 258
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 867
$ws.Range("F4").Value = 245
$ws.Range("F6").Value = 31
$ws.Range("F7").Value = 45
$ws.Range("F8").Value = 3
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 867
$ws.Range("F7").Value = 245
$ws.Range("F11").Value = 31
$ws.Range("F12").Value = 45
$ws.Range("F14").Value = 16
$ws.Range("F21").Value = 2675
$ws.Range("F23").Value = 1754
$ws.Range("F24").Value = 137
$ws.Range("F25").Value = 607
$ws.Range("F26").Value = 273
$ws.Range("F27").Value = 686
$ws.Range("F28").Value = 5030
$ws.Range("F31").Value = 3385
$ws.Range("F35").Value = 37
$ws.Range("F36").Value = 2417
$ws.Range("F38").Value = 366
$ws.Range("F41").Value = 480
$ws.Range("F42").Value = 1295
$ws.Range("F43").Value = 139
$ws.Range("F44").Value = 258
$ws.Range("F49").Value = 54
$ws.Range("F50").Value = 1430
